# TC06_Canine_Filter_Breed-BelgMalin.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The CasesTab query (row 2, column B) incorrectly included a trailing
# clause that joined to a (co:cohort) node and returned a `Cohort` column
# that isn't part of this query's RETURN contract. Remove that clause so
# the query text matches the corrected version.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$fixedCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Belgian Malinois']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $fixedCasesQuery

# Reflect the author's final view state: scrolled back to the top,
# zoomed in to 115%, with B2 (the row they just fixed) selected.
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 115
